# "In Class Demonstrations" table update:
#   - insert a new "$\pi$" / "$\pi$ Lag" column between the existing U and FFR
#     columns (and corresponding row)
#   - drop the old "Constant" and "r2_adj" rows
#   - refresh all the regression coefficients with the new run's values

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the trailing "r2_adj" row - it's no longer part of the table.
$ws.Rows("5").Delete()

# Make room for the new "$\pi$" column between U (B) and FFR (old C, now D).
$ws.Columns("C").Insert()

# --- Header row -------------------------------------------------------
$ws.Range("C1").Value = "$\pi$"

# --- Row labels ---------------------------------------------------------
# Row 3 used to be "FFR Lag"; it now holds the new "$\pi$ Lag" data.
$ws.Range("A3").Value = "$\pi$ Lag"
# Row 4 used to be "Constant"; it now holds the "FFR Lag" data.
$ws.Range("A4").Value = "FFR Lag"

# A handful of the new coefficients (0.43, -0.507, 1.049, 0.011) are plain
# numerals with no trailing "*"/"**"/"***" significance markers, so Excel
# would otherwise read them in as numbers. They need to stay text (matching
# how every other coefficient in this table is stored), so they're entered
# with a leading apostrophe and then reset to the plain/default style so no
# left-over "stored as text" formatting sticks to the cell.

# --- U column (B) --------------------------------------------------------
$ws.Range("B2").Value = "'0.43"
$ws.Range("B2").Style = $ws.Range("A1").Style
$ws.Range("B3").Value = "-0.558**"
$ws.Range("B4").Value = "-0.186**"

# --- $\pi$ column (C, newly inserted) ------------------------------------
$ws.Range("C2").Value = "'-0.507"
$ws.Range("C2").Style = $ws.Range("A1").Style
$ws.Range("C3").Value = "-0.993***"
$ws.Range("C4").Value = "'0.011"
$ws.Range("C4").Style = $ws.Range("A1").Style

# --- FFR column (D) --------------------------------------------------------
$ws.Range("D2").Value = "-2.552*"
$ws.Range("D3").Value = "'1.049"
$ws.Range("D3").Style = $ws.Range("A1").Style
$ws.Range("D4").Value = "0.475*"
